$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newMain = '![main_banner not_rounded](data/img/main_banner{dark_mode}.png?v=1)

# datannur, le catalogue de données portable

Permet de centraliser, rechercher et visualiser les informations sur une collection de jeux de données

Pour améliorer l’organisation des données et faciliter leur partage et leur documentation

Simple et flexible, s’intègre rapidement dans tous types d’environnement


- **Facile** :
Aucune installation ou configuration nécessaire. datannur est le catalogue le plus simple à implémenter et maintenir

- **Portable** :
Fonctionne partout (local, cloud, disque partagé), un dossier que l’on peut copier, déplacer, envoyer et ouvrir avec n’importe quel navigateur

- **Complet** :
Flexible, complet et structuré autour de 7 concepts avec un niveau de détail important : Institution, Dossier, Mot clé, Doc, Dataset, Variable et Modalité

- **Sécurisé** :
Parce qu’elle est une simple interface HTML isolée dans le navigateur, l’application ne peut rien modifier sur la machine et ne pose ainsi aucun risque

Pour davantage d''information : [datannur.com](https://datannur.com)

La version ici présente est un prototype en cours de développement et d''expérimentation. Les données utilisées sont fictives et uniquement à usage de test et de développement. Question ou suggestion : [contact@datannur.com](mailto:contact@datannur.com).'

$ws.Range("B3").Value = $newMain
[void]$ws.Range("B3").Select()
